$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: sdmx/iaest-measure semantic mapping per column
$ws.Range("A2").Value = "null"
$ws.Range("B2").Value = "sdmx-dimension:refArea"
$ws.Range("C2").Value = "sdmx-dimension:refArea"
$ws.Range("D2").Value = "iaest-measure:porcentaje"
$ws.Range("E2").Value = "sdmx-dimension:refArea"
$ws.Range("F2").Value = "iaest-measure:sexo"
$ws.Range("G2").Value = "null"
$ws.Range("H2").Value = "iaest-measure:personas"
$ws.Range("I2").Value = "iaest-measure:edad-grupos-quinquenales"
$ws.Range("J2").Value = "sdmx-dimension:refArea"
$ws.Range("K2").Value = "null"

# Row 3: dim/medida classification per column
$ws.Range("A3").Value = "null"
$ws.Range("B3").Value = "dim"
$ws.Range("C3").Value = "dim"
$ws.Range("D3").Value = "medida"
$ws.Range("E3").Value = "dim"
$ws.Range("F3").Value = "medida"
$ws.Range("G3").Value = "null"
$ws.Range("H3").Value = "medida"
$ws.Range("I3").Value = "medida"
$ws.Range("J3").Value = "dim"
$ws.Range("K3").Value = "null"

# Row 4: URI / xsd type per column
$ws.Range("A4").Value = "null"
$ws.Range("B4").Value = "URI-Municipio"
$ws.Range("C4").Value = "URI-Provincia"
$ws.Range("D4").Value = "xsd:int"
$ws.Range("E4").Value = "URI-Comunidad"
$ws.Range("F4").Value = "xsd:int"
$ws.Range("G4").Value = "null"
$ws.Range("H4").Value = "xsd:int"
$ws.Range("I4").Value = "xsd:int"
$ws.Range("J4").Value = "URI-comarca"
$ws.Range("K4").Value = "null"

# Row 5 no longer exists in the curated dataset - remove it entirely
$ws.Range("A5:K5").EntireRow.Delete()
